# Replace every occurrence of the variable name "congenital" with
# "misc_long_term" across all worksheets in the workbook.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count

    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Text -eq "congenital") {
                $cell.Value = "misc_long_term"
            }
        }
    }
}
